# Updated C3DC phs000466 queries
#
# The TreatmentTab query (cell B5 on Sheet1) wrapped the REPLACE() call in an
# unnecessary CONCAT(); simplify it to a plain REPLACE() call.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$cell = $ws.Range("B5")
$old = $cell.Value()
$new = $old.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent""", "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent""")
$cell.Value = $new

# Leave the selection on B2, matching the saved workbook view state.
$ws.Range("B2").Select() | Out-Null
